$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH347-1"
$ws.Range("C2").Value = "14/08/2006 FROM MANKI & CAROLL... ( TO BE SORTED)"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2006"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 33K | GRAP COUNT NUMER: NONE"

$ws.Range("H20").Select()
